# Updated cryptos list on Wed Mar 27 21:26:53 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "1.00", "64.90",
# "3.517.50" using "." as a thousands separator). Force the cell to Text
# format before assigning so Excel keeps the exact printed digits instead
# of silently re-parsing the string as a number and dropping/altering it.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.228.60'
$ws.Range('E2').Value = '  -0.95%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.517.50'
$ws.Range('E3').Value = '  -1.68%  '

$ws.Range('E4').Value = '  -0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '573.08'
$ws.Range('E5').Value = '  -0.91%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '185.38'
$ws.Range('E6').Value = '  -2.66%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.509.87'
$ws.Range('E7').Value = '  -1.84%  '

$ws.Range('E8').Value = '  -2.80%  '

$ws.Range('E9').Value = '  +0.06%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.185'
$ws.Range('E10').Value = '  +2.20%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.657'
$ws.Range('E11').Value = '  -0.87%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.37'
$ws.Range('E12').Value = '  -2.58%  '

$ws.Range('E13').Value = '  -1.37%  '

$ws.Range('E14').Value = '  -1.75%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.079.42'
$ws.Range('E15').Value = '  -1.85%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '19.48'
$ws.Range('E16').Value = '  -2.34%  '

$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.517.98'
$ws.Range('E17').Value = '  -1.66%  '

$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.138.83'
$ws.Range('E18').Value = '  -1.08%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.35'
$ws.Range('E19').Value = '  -2.60%  '

$ws.Range('E20').Value = '  -1.05%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '545.81'
$ws.Range('E21').Value = '  +14.01%  '

$ws.Range('E22').Value = '  -2.68%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '18.72'
$ws.Range('E23').Value = '  -3.28%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.99'
$ws.Range('E24').Value = '  -1.09%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.44'
$ws.Range('E25').Value = '  +1.29%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '94.42'
$ws.Range('E26').Value = '  -1.78%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.92'
$ws.Range('E27').Value = '  -2.54%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.89'
$ws.Range('E28').Value = '  -1.53%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.14'
$ws.Range('E29').Value = '  -2.87%  '

$ws.Range('E30').Value = '  -1.38%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.25'
$ws.Range('E31').Value = '  -5.60%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.66'
$ws.Range('E32').Value = '  +3.49%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '64.90'
$ws.Range('E33').Value = '  -1.99%  '

$ws.Range('E34').Value = '  -3.93%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '558.37'
$ws.Range('E35').Value = '  -4.20%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '38.10'
$ws.Range('E36').Value = '  -2.24%  '

$ws.Range('B37').Value = 'Fetch.AI'
$ws.Range('C37').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.07'
$ws.Range('E37').Value = '  +7.50%  '

$ws.Range('B38').Value = 'TheGraph'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.401'
$ws.Range('E38').Value = '  +1.20%  '

$ws.Range('B39').Value = 'Dai'
$ws.Range('C39').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.00'
$ws.Range('E39').Value = '  +0.08%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0₃0767'
$ws.Range('E40').Value = '  -4.40%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.12'
$ws.Range('E41').Value = '  -3.85%  '

$ws.Range('E42').Value = '  -2.71%  '

$ws.Range('E43').Value = '  -3.48%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.283.16'
$ws.Range('E44').Value = '  +1.64%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.55'
$ws.Range('E45').Value = '  +6.01%  '

$ws.Range('E46').Value = '  -2.77%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0443'
$ws.Range('E47').Value = '  -0.10%  '

$ws.Range('E48').Value = '  -2.25%  '

$ws.Range('E49').Value = '  -4.38%  '

$ws.Range('E50').Value = '  -0.14%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '138.59'
$ws.Range('E51').Value = '  +3.30%  '
